# Sale fixes 1 (#324): add a new "Seller Signatory Emails" column (I) to the
# Offers sheet, populated with per-row signer emails (with hyperlinks for the
# rows that already used a mailto: hyperlink pattern elsewhere on the row),
# and tidy up a redundant cell style that Excel had been carrying around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column at I. Everything that used to be columns I..O
#    (Bank Account .. Client Id) shifts right to J..P.
# ---------------------------------------------------------------------
$ws.Columns("I:I").Insert()

# Give the new column roughly the same width as its neighbour (H) instead
# of the sheet default.
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth

# ---------------------------------------------------------------------
# 2. Header + values for the new "Seller Signatory Emails" column.
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "Seller Signatory Emails"
$ws.Range("I2").Value = "emp1@myfirm.com"
$ws.Range("I3").Value = "emp2@myfirm.com"
$ws.Range("I4").Value = "emp3@myfirm.com"
$ws.Range("I5").Value = "emp4@myfirm.com"
$ws.Range("I6").Value = "emp1@investor1.com"
$ws.Range("I7").Value = "emp1@investor2.com"

# ---------------------------------------------------------------------
# 3. Hyperlink the three addresses that already have a mailto: hyperlink
#    elsewhere on their row (same targets used by column D's links).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:emp3@myfirm.com")
$ws.Hyperlinks.Add($ws.Range("I6"), "mailto:emp1@investor1.com")
$ws.Hyperlinks.Add($ws.Range("I7"), "mailto:emp1@investor2.com")

# ---------------------------------------------------------------------
# 4. Match the new column's look to column D (same content pattern:
#    plain email text, or a hyperlinked email). Re-applied last so it
#    wins over whatever default formatting the hyperlink step above
#    introduced.
# ---------------------------------------------------------------------
$ws.Range("D2:D7").Copy()
$ws.Range("I2:I7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. The old "Update Only" column (now N) -- and the header cells of the
#    two columns after it (now O1, P1) -- had been carrying a redundant
#    duplicate cell style (identical apart from a no-op applyFill flag).
#    Re-apply the plain style so they match the rest of the header row.
# ---------------------------------------------------------------------
$ws.Range("N1:N7").Style = $ws.Range("J1").Style
$ws.Range("O1").Style = $ws.Range("J1").Style
$ws.Range("P1").Style = $ws.Range("J1").Style
